# Auto-generated edit script: updates cached market-price / profit
# columns (H:N) on each leve sheet to match the refreshed data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 266.1111
$ws.Range("I4").Value = 261.875
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 261.875
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -147.875
$ws.Range("N4").Value = -528
$ws.Range("H19").Value = 400
$ws.Range("J19").Value = 400
$ws.Range("L19").Value = 400
$ws.Range("N19").Value = -750
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532
$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766
$ws.Range("H26").Value = 664.5
$ws.Range("I26").Value = 664.5
$ws.Range("K26").Value = 664.5
$ws.Range("M26").Value = -320.5
$ws.Range("H29").Value = 2291.125
$ws.Range("I29").Value = 194.33333
$ws.Range("K29").Value = 582.99999
$ws.Range("M29").Value = -301.99999
$ws.Range("H31").Value = 11.5
$ws.Range("I31").Value = 11.5
$ws.Range("K31").Value = 34.5
$ws.Range("M31").Value = 195.5
$ws.Range("H98").Value = 3576.4443
$ws.Range("I98").Value = 3455.5715
$ws.Range("K98").Value = 3455.5715
$ws.Range("M98").Value = -1957.5715
$ws.Range("H101").Value = 1079.2222
$ws.Range("I101").Value = 591.5
$ws.Range("K101").Value = 1774.5
$ws.Range("M101").Value = -152.5
$ws.Range("H122").Value = 3576.4443
$ws.Range("I122").Value = 3455.5715
$ws.Range("K122").Value = 10366.7145
$ws.Range("M122").Value = -7916.7145
$ws.Range("H132").Value = 4362.4
$ws.Range("J132").Value = 3961.5
$ws.Range("L132").Value = 11884.5
$ws.Range("N132").Value = -16944.5
$ws.Range("I141").Value = 4301.1055
$ws.Range("J141").Value = 9997
$ws.Range("K141").Value = 12903.3165
$ws.Range("L141").Value = 29991
$ws.Range("M141").Value = -7723.316499999999
$ws.Range("N141").Value = -40351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1899.3334
$ws.Range("I2").Value = 2674
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 2674
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -2561
$ws.Range("N2").Value = -576
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -623
$ws.Range("H74").Value = 4600
$ws.Range("I74").Value = 4600
$ws.Range("K74").Value = 4600
$ws.Range("M74").Value = -3726
$ws.Range("H77").Value = 4600
$ws.Range("I77").Value = 4600
$ws.Range("K77").Value = 23000
$ws.Range("M77").Value = -18632
$ws.Range("H116").Value = 1899.3334
$ws.Range("I116").Value = 2674
$ws.Range("J116").Value = 350
$ws.Range("K116").Value = 2674
$ws.Range("L116").Value = 350
$ws.Range("M116").Value = -380
$ws.Range("N116").Value = -4938
$ws.Range("H122").Value = 2485.125
$ws.Range("I122").Value = 2426
$ws.Range("K122").Value = 7278
$ws.Range("M122").Value = -4828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1899.3334
$ws.Range("I3").Value = 2674
$ws.Range("J3").Value = 350
$ws.Range("K3").Value = 2674
$ws.Range("L3").Value = 350
$ws.Range("M3").Value = -2560
$ws.Range("N3").Value = -578
$ws.Range("H86").Value = 3719.9
$ws.Range("J86").Value = 3859.8
$ws.Range("L86").Value = 3859.8
$ws.Range("N86").Value = -6105.8
$ws.Range("H89").Value = 3719.9
$ws.Range("J89").Value = 3859.8
$ws.Range("L89").Value = 19299
$ws.Range("N89").Value = -30531
$ws.Range("H134").Value = 5486.9165
$ws.Range("I134").Value = 6349.222
$ws.Range("K134").Value = 19047.666
$ws.Range("M134").Value = -16512.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4197.722
$ws.Range("J16").Value = 6855.75
$ws.Range("L16").Value = 6855.75
$ws.Range("N16").Value = -7429.75
$ws.Range("H58").Value = 3632.5454
$ws.Range("I58").Value = 3196.8
$ws.Range("J58").Value = 3995.6667
$ws.Range("K58").Value = 3196.8
$ws.Range("L58").Value = 3995.6667
$ws.Range("M58").Value = -2993.8
$ws.Range("N58").Value = -4401.6667
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H105").Value = 1004.6667
$ws.Range("I105").Value = 1004.6667
$ws.Range("K105").Value = 1004.6667
$ws.Range("M105").Value = 742.3333
$ws.Range("H113").Value = 4197.722
$ws.Range("J113").Value = 6855.75
$ws.Range("L113").Value = 6855.75
$ws.Range("N113").Value = -11195.75
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H135").Value = 30000
$ws.Range("I135").Value = 30000
$ws.Range("K135").Value = 30000
$ws.Range("M135").Value = -24930
$ws.Range("H136").Value = 3632.5454
$ws.Range("I136").Value = 3196.8
$ws.Range("J136").Value = 3995.6667
$ws.Range("K136").Value = 9590.400000000001
$ws.Range("L136").Value = 11987.0001
$ws.Range("M136").Value = -7040.400000000001
$ws.Range("N136").Value = -17087.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3898.5
$ws.Range("J22").Value = 3898.5
$ws.Range("L22").Value = 11695.5
$ws.Range("N22").Value = -12033.5
$ws.Range("H27").Value = 3898.5
$ws.Range("J27").Value = 3898.5
$ws.Range("L27").Value = 11695.5
$ws.Range("N27").Value = -11899.5
$ws.Range("H39").Value = 841.6667
$ws.Range("I39").Value = 841.6667
$ws.Range("K39").Value = 2525.0001
$ws.Range("M39").Value = -2231.0001
$ws.Range("H41").Value = 1333.3334
$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2909
$ws.Range("H51").Value = 105
$ws.Range("I51").Value = 105
$ws.Range("K51").Value = 315
$ws.Range("M51").Value = 145
$ws.Range("H57").Value = 4999.5
$ws.Range("I57").Value = 4999.5
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 14998.5
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0
$ws.Range("M57").Value = -14439.5
$ws.Range("H113").Value = 934.9231
$ws.Range("J113").Value = 1137.3334
$ws.Range("L113").Value = 3412.0002
$ws.Range("N113").Value = -7752.0002
$ws.Range("H117").Value = 841.7
$ws.Range("J117").Value = 1187.25
$ws.Range("L117").Value = 3561.75
$ws.Range("N117").Value = -10445.75
$ws.Range("H128").Value = 339999
$ws.Range("I128").Value = 339999
$ws.Range("K128").Value = 1019997
$ws.Range("M128").Value = -1015017

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 105742
$ws.Range("I22").Value = 1475
$ws.Range("K22").Value = 1475
$ws.Range("M22").Value = -946
$ws.Range("H126").Value = 2399.5
$ws.Range("I126").Value = 2399.6667
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 7199.000100000001
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -4729.000100000001
$ws.Range("N126").Value = -12137
$ws.Range("H132").Value = 1234
$ws.Range("I132").Value = 1234
$ws.Range("K132").Value = 3702
$ws.Range("M132").Value = -1172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 823.5833
$ws.Range("I22").Value = 898.36365
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = 898.36365
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = -603.36365
$ws.Range("N22").Value = -591
$ws.Range("H27").Value = 823.5833
$ws.Range("I27").Value = 898.36365
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 898.36365
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = -791.36365
$ws.Range("N27").Value = -215
$ws.Range("H46").Value = 3537.4614
$ws.Range("I46").Value = 2498
$ws.Range("K46").Value = 2498
$ws.Range("M46").Value = -2310
$ws.Range("H61").Value = 774.5
$ws.Range("I61").Value = 774.5
$ws.Range("K61").Value = 774.5
$ws.Range("M61").Value = -572.5
$ws.Range("H113").Value = 774.5
$ws.Range("I113").Value = 774.5
$ws.Range("K113").Value = 774.5
$ws.Range("M113").Value = 1395.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 389.58334
$ws.Range("I81").Value = 370.45456
$ws.Range("K81").Value = 740.90912
$ws.Range("M81").Value = 320.09088
$ws.Range("H84").Value = 389.58334
$ws.Range("I84").Value = 370.45456
$ws.Range("K84").Value = 3704.5456
$ws.Range("M84").Value = 1599.4544
$ws.Range("H113").Value = 35134.668
$ws.Range("I113").Value = 51350.5
$ws.Range("J113").Value = 2703
$ws.Range("K113").Value = 154051.5
$ws.Range("L113").Value = 8109
$ws.Range("M113").Value = -151881.5
$ws.Range("N113").Value = -12449
$ws.Range("H126").Value = 4995.25
$ws.Range("I126").Value = 4995.25
$ws.Range("K126").Value = 14985.75
$ws.Range("M126").Value = -12515.75
